$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# Move the existing "Debug Mode / 로딩 화면" rounded rectangle shape
$s.Shapes.Item(1).Left = 195943
$s.Shapes.Item(1).Top = 114300

# Add a new rounded rectangle shape for "운동 전 안내 화면"
$newShape = $s.Shapes.AddShape(5, 5954486, 57150, 5595257, 6743700)
$newShape.Name = "사각형: 둥근 모서리 1"
$newShape.TextFrame.TextRange.Text = "운동 전 안내 화면"
$newShape.TextFrame.TextRange.Font.Size = 36
